# Updated symbol list on Sun Jan 29 22:06:42 UTC 2023 with GitHub Actions
# Refresh the crypto price/volume/hour columns for each row of the table
# (values are written with a leading apostrophe so Excel keeps them as
#  literal text, matching the inlineStr cells already used in the sheet)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'317.86"
$ws.Range("E2").Value = "'3.99%"
$ws.Range("G2").Value = "'22"

# Row 3
$ws.Range("D3").Value = "'39.66"
$ws.Range("E3").Value = "'2.22%"
$ws.Range("G3").Value = "'22"

# Row 4
$ws.Range("D4").Value = "'5.137"
$ws.Range("E4").Value = "'0.64%"
$ws.Range("G4").Value = "'22"

# Row 5
$ws.Range("D5").Value = "'0.08210"
$ws.Range("E5").Value = "'1.72%"
$ws.Range("G5").Value = "'22"

# Row 6
$ws.Range("D6").Value = "'2.048"
$ws.Range("E6").Value = "'6.04%"
$ws.Range("G6").Value = "'22"

# Row 7
$ws.Range("D7").Value = "'8.354"
$ws.Range("E7").Value = "'4.27%"
$ws.Range("G7").Value = "'22"

# Row 8
$ws.Range("D8").Value = "'4.315"
$ws.Range("E8").Value = "'2.58%"
$ws.Range("G8").Value = "'22"

# Row 9
$ws.Range("D9").Value = "'0.9404"
$ws.Range("E9").Value = "'1.41%"
$ws.Range("G9").Value = "'22"

# Row 10
$ws.Range("D10").Value = "'0.1359"
$ws.Range("E10").Value = "'-5.33%"
$ws.Range("G10").Value = "'22"

# Row 11
$ws.Range("D11").Value = "'0.1994"
$ws.Range("E11").Value = "'4.15%"
$ws.Range("G11").Value = "'22"

# Row 12
$ws.Range("E12").Value = "'0.50%"
$ws.Range("G12").Value = "'22"

# Row 13
$ws.Range("D13").Value = "'0.03511"
$ws.Range("E13").Value = "'0.15%"
$ws.Range("G13").Value = "'22"

# Row 14
$ws.Range("D14").Value = "'0.09798"
$ws.Range("E14").Value = "'0.10%"
$ws.Range("G14").Value = "'22"

# Row 15
$ws.Range("D15").Value = "'0.001411"
$ws.Range("E15").Value = "'1.07%"
$ws.Range("G15").Value = "'22"

# Row 16
$ws.Range("D16").Value = "'0.006239"
$ws.Range("E16").Value = "'4.96%"
$ws.Range("G16").Value = "'22"

# Row 17
$ws.Range("D17").Value = "'3.682"
$ws.Range("E17").Value = "'-2.65%"
$ws.Range("G17").Value = "'22"

# Row 18
$ws.Range("G18").Value = "'22"

# Row 19
$ws.Range("D19").Value = "'0.3479"
$ws.Range("E19").Value = "'0.51%"
$ws.Range("G19").Value = "'22"

# Row 20
$ws.Range("D20").Value = "'0.1321"
$ws.Range("E20").Value = "'-0.49%"
$ws.Range("G20").Value = "'22"

# Row 21
$ws.Range("D21").Value = "'4.989"
$ws.Range("E21").Value = "'6.50%"
$ws.Range("G21").Value = "'22"

# Row 22
$ws.Range("D22").Value = "'0.2450"
$ws.Range("E22").Value = "'1.36%"
$ws.Range("G22").Value = "'22"

# Row 23
$ws.Range("D23").Value = "'0.04355"
$ws.Range("E23").Value = "'-0.50%"
$ws.Range("G23").Value = "'22"

# Row 24
$ws.Range("E24").Value = "'0.31%"
$ws.Range("G24").Value = "'22"

# Row 25
$ws.Range("D25").Value = "'0.004799"
$ws.Range("E25").Value = "'12.28%"
$ws.Range("G25").Value = "'22"

# Row 26
$ws.Range("E26").Value = "'-0.11%"
$ws.Range("G26").Value = "'22"

# Row 27
$ws.Range("D27").Value = "'0.0003998"
$ws.Range("E27").Value = "'-10.10%"
$ws.Range("G27").Value = "'22"

# Row 28
$ws.Range("G28").Value = "'22"

# Row 29
$ws.Range("G29").Value = "'22"

# Row 30
$ws.Range("G30").Value = "'22"

# Row 31
$ws.Range("G31").Value = "'22"

# Row 32
$ws.Range("G32").Value = "'22"

# Row 33
$ws.Range("G33").Value = "'22"

# Row 34
$ws.Range("G34").Value = "'22"

# Row 35
$ws.Range("G35").Value = "'22"

# Row 36
$ws.Range("G36").Value = "'22"

# Row 37
$ws.Range("G37").Value = "'22"

# Row 38
$ws.Range("G38").Value = "'22"

# Row 39
$ws.Range("D39").Value = "'0.02258"
$ws.Range("E39").Value = "'11.11%"
$ws.Range("G39").Value = "'22"

# Row 40
$ws.Range("D40").Value = "'0.05190"
$ws.Range("E40").Value = "'2.88%"
$ws.Range("G40").Value = "'22"

# Row 41
$ws.Range("D41").Value = "'0.007765"
$ws.Range("E41").Value = "'3.22%"
$ws.Range("G41").Value = "'22"

# Row 42
$ws.Range("D42").Value = "'0.009871"
$ws.Range("E42").Value = "'1.46%"
$ws.Range("G42").Value = "'22"

# Row 43
$ws.Range("D43").Value = "'0.1405"
$ws.Range("E43").Value = "'4.89%"
$ws.Range("G43").Value = "'22"

# Row 44
$ws.Range("D44").Value = "'0.002081"
$ws.Range("E44").Value = "'-1.13%"
$ws.Range("G44").Value = "'22"

# Row 45
$ws.Range("D45").Value = "'0.009656"
$ws.Range("E45").Value = "'-2.47%"
$ws.Range("G45").Value = "'22"

# Row 46
$ws.Range("D46").Value = "'0.00006597"
$ws.Range("E46").Value = "'6.15%"
$ws.Range("G46").Value = "'22"

# Row 47
$ws.Range("G47").Value = "'22"

# Row 48
$ws.Range("E48").Value = "'2.41%"
$ws.Range("G48").Value = "'22"

# Row 49
$ws.Range("D49").Value = "'0.001689"
$ws.Range("G49").Value = "'22"

# Row 50
$ws.Range("G50").Value = "'22"

# Row 51
$ws.Range("G51").Value = "'22"
